# Settings back to qfin22 presentation - recalibration
# Updates calibration parameter values (and their associated autocorrelation
# arrays stored as text) across the AR, SETAR, GARCH, TARCH and AR_TARCH sheets.

$wb = $excel.ActiveWorkbook

# ----- AR sheet -----
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.005754839538796336
$ws.Range("B3").Value = 0.795922957307973
$ws.Range("B4").Value = 3.347412037893251
$ws.Range("B5").Value = "[0.9999999999999999, 0.27650242874803493, 0.2473059352921494, 0.22368876536446608, 0.2650651106805935, 0.4148231065895015, 0.23266783219559437, 0.24716826330919148, 0.19328402907279793, 0.21036520904878678, 0.17430397655654778, 0.19441324726301387, 0.18868717834654988, 0.1990158199579622, 0.1956654914708275, 0.2097120307742413, 0.1855564931629684, 0.19637246978488584, 0.17975199203381023, 0.17428121219450082]"

# ----- SETAR sheet -----
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.5885009830465028
$ws.Range("B5").Value = 0.7262573834541782
$ws.Range("B6").Value = 2.530991530374362
$ws.Range("B7").Value = 0.5871451601940401
$ws.Range("B8").Value = 0.7156586940581252
$ws.Range("B9").Value = 1.992191975033456
$ws.Range("B10").Value = "[1.0, 0.2865652493659193, 0.2536966090007755, 0.23647302114734567, 0.26388517285298846, 0.28534726148662526, 0.2126390971948123, 0.22206088841954952, 0.18710997547316743, 0.17013249022154797, 0.1504470990959502, 0.18788480141933467, 0.1667622997722943, 0.20646199863349443, 0.16459115766088608, 0.18758420582371566, 0.1713078106224247, 0.17614104180344584, 0.1544029417985861, 0.15730089318608456]"

# ----- GARCH sheet -----
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.008841838364523173
$ws.Range("B3").Value = 0.01889802222074789
$ws.Range("B4").Value = 0.0792058023404483
$ws.Range("B5").Value = 0.9194519060039666
$ws.Range("B6").Value = "[1.0, 0.013033516628285575, -0.023964128964601913, -0.05253996601624041, -0.03387497693178653, 0.27604735862502033, 0.015955780864374405, 0.01821931347351153, -0.030449947667353167, -0.027320762833466457, -0.0699261422457061, -0.04847873239455806, -0.03288594114235811, -0.020915784563769493, 0.030542870029760333, 0.00495460224119161, -0.009515930843045126, -0.04552672124974095, -0.0317197586874379, -0.021228542782399905]"

# ----- TARCH sheet -----
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.01354273611591852
$ws.Range("B3").Value = 0.01912439436641028
$ws.Range("B4").Value = 0.07201274384547182
$ws.Range("B5").Value = 0.92041911077476
$ws.Range("B6").Value = "[1.0, 0.013934789526919726, -0.023613640150222596, -0.051903422328022225, -0.03325942608738673, 0.275855881748233, 0.01663166510376376, 0.01853587090887995, -0.030261012496108363, -0.026995927541476665, -0.06947620845214035, -0.04816459290216756, -0.03257512017260114, -0.02070943480622141, 0.031061685139023108, 0.0053927829773057195, -0.009149522523800737, -0.04554608681735547, -0.031238995389322486, -0.02119675418038461]"
$ws.Range("B7").Value = 0.01192603574287862

# ----- AR_TARCH sheet -----
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.003335417073743878
$ws.Range("B3").Value = 0.01334447024231665
$ws.Range("B4").Value = 0.08809867016160182
$ws.Range("B5").Value = 0.922332784293413
$ws.Range("B6").Value = "[1.0, 0.03271258947821372, 0.0012238041781187107, -0.023441013082385363, 0.005368262609537936, 0.20493357806707266, 0.0027742686247121968, 0.01200685294674964, -0.03481047948004237, -0.036417019883907956, -0.05873577257063197, -0.04297618931346965, -0.04466499820083468, -0.03604702878248634, 0.02029002251354862, -0.009950318283285764, -0.019808412445069697, -0.036840995848426726, -0.023401635521485505, -0.020653724902000105]"
$ws.Range("B7").Value = -0.02086290891002961
$ws.Range("B9").Value = 0.7899467740281839
